$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange

# Order matters: do the size replacement first (S30 -> S31) since "S30" does not
# collide with any distance token, then the distance replacements. Each call only
# touches cells containing the literal old token, and only the matched substring
# is replaced (xlPart), preserving the rest of composite strings such as
# "Face08_D64_S30" -> "Face08_D69_S31" and "Face08_D64_S30_l.png" -> "Face08_D69_S31_l.png".

$used.Replace("S30", "S31", [Microsoft.Office.Interop.Excel.XlLookAt]::xlPart)
$used.Replace("D64", "D69", [Microsoft.Office.Interop.Excel.XlLookAt]::xlPart)
$used.Replace("D51", "D55", [Microsoft.Office.Interop.Excel.XlLookAt]::xlPart)
$used.Replace("D80", "D86", [Microsoft.Office.Interop.Excel.XlLookAt]::xlPart)
